$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows before row 356. This shifts the existing rows 356-374
# down to 358-376 (values, formats and styles move with them), matching the
# new dimension A1:R376 and the fact that former rows 356-374 reappear,
# unchanged, two rows further down (with the former 373/374 content also
# duplicated onto the brand-new trailing rows 375/376).
$ws.Rows("356:357").Insert()

# The Insert() above left rows 356-357 completely blank. Fill them in with
# a new week's worth of data for the same market/product (same as the
# unchanged columns of the rows that used to sit there, now at 358/359,
# plus the new date, volume, price and $/Kg values from the diff).
$ws.Range("A356").Value = 1
$ws.Range("B356").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C356").Value = "Arica y Parinacota"
$ws.Range("D356").Value = 44746
$ws.Range("E356").Value = 15
$ws.Range("F356").Value = 100112032
$ws.Range("G356").Value = "Zapallo italiano"
$ws.Range("H356").Value = "Huracán"
$ws.Range("I356").Value = "Primera"
$ws.Range("J356").Value = 150
$ws.Range("K356").Value = 8000
$ws.Range("L356").Value = 8500
$ws.Range("M356").Value = 8250
$ws.Range("N356").Value = "$/caja 70 unidades"
$ws.Range("O356").Value = "Región de Arica y Parinacota"
$ws.Range("P356").Value = 118
$ws.Range("Q356").Value = 70
$ws.Range("R356").Value = "Hortaliza"

$ws.Range("A357").Value = 1
$ws.Range("B357").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C357").Value = "Arica y Parinacota"
$ws.Range("D357").Value = 44746
$ws.Range("E357").Value = 15
$ws.Range("F357").Value = 100112032
$ws.Range("G357").Value = "Zapallo italiano"
$ws.Range("H357").Value = "Huracán"
$ws.Range("I357").Value = "Segunda"
$ws.Range("J357").Value = 150
$ws.Range("K357").Value = 7000
$ws.Range("L357").Value = 7500
$ws.Range("M357").Value = 7250
$ws.Range("N357").Value = "$/caja 100 unidades"
$ws.Range("O357").Value = "Región de Arica y Parinacota"
$ws.Range("P357").Value = 72
$ws.Range("Q357").Value = 100
$ws.Range("R357").Value = "Hortaliza"
